# Update the second row ("best accuracy" series) on the active sheet with
# newly trained values (HCN02 run). Only the B2:U2 numeric cells change;
# everything else (headers, layout, styles) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60.36324782503976
$ws.Range("C2").Value = 50.94818373521169
$ws.Range("D2").Value = 46.62126070923276
$ws.Range("E2").Value = 45.69978614648183
$ws.Range("F2").Value = 44.23076907793681
$ws.Range("G2").Value = 43.70993574460348
$ws.Range("H2").Value = 43.01549130015903
$ws.Range("I2").Value = 43.01549130015903
$ws.Range("J2").Value = 42.58814089828067
$ws.Range("K2").Value = 42.58814089828067
$ws.Range("L2").Value = 42.41452978716956
$ws.Range("M2").Value = 42.41452978716956
$ws.Range("N2").Value = 42.41452978716956
$ws.Range("O2").Value = 42.41452978716956
$ws.Range("P2").Value = 42.41452978716956
$ws.Range("Q2").Value = 42.41452978716956
$ws.Range("R2").Value = 42.41452978716956
$ws.Range("S2").Value = 42.41452978716956
$ws.Range("T2").Value = 42.41452978716956
$ws.Range("U2").Value = 42.41452978716956
